# Append row 47 to the daily trading_data series (date serial 43857,
# i.e. 2020-01-27 -- the next trading day after row 46's 2020-01-24) across
# the full A:RR column range (486 market-value columns), and clean up the
# stray negative-zero left in MI46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- normalize the stray "-0" left over in the previous row ---
$ws.Range("MI46").Value = 0

# --- give A47 the same date style as A46 (centered/bordered, custom
#     "YYYY-MM-DD HH:MM:SS" number format) by copying formats only ---
$ws.Range("A46").Copy()
$ws.Range("A47").PasteSpecial(-4122)

# --- the 486 values for row 47, in column order A, B, C, ... RR ---
$rowValues = @(
    43857,0,141.5619178077402,0,0,243.5430523819751,0,128.7960589909675,1160.643911413938,49.56160694437676,903.0865242788241,0,0,3.407910497242101,607.3636272486135,
    825.4785559185738,603.308850223541,27.19358954552806,0,0,382.0965802827436,0,0,952.4094435083773,229.8848220045544,36.25957421451911,38.35317690150693,0,0,4011.895915323778,
    328.4823254914574,1778.212969740873,559.2897039746528,284.3219657333249,0,40.75636455702903,420.4597212613269,319.507022943754,0,0,0,0,483.3175565002795,194.4662418831975,0,
    264.2273604914741,0,0,416.9161187198258,185.077923127406,731.7208976805841,541.8999834611373,0,59.03506789851102,213.7344263405161,493.6214727038487,0,434.5621599604853,60.66777948902764,141.0533882936172,
    175.5113662139754,46.53918375446253,0,116.5022293858265,1549.481234404241,588.6857235636483,0,766.2275865044354,68.18279498210079,141.6635448428823,154.2982081822491,0,0,108.4464659463161,46.01032394546784,
    2.08771481674073,67.24685074317733,231.1734089319239,0,0,184.6613057117234,0,702.3177803970052,0,0,729.5236481995167,908.091654480643,348.3210565931838,0,450.4286800693229,
    581.3479624683254,0,0,475.2822770057195,0,381.1003908938076,0,908.3968043428185,0,534.8377662120251,0,0,0,105.4465455642277,0,
    53.54824796936703,0,0,0,0,0,173.847253633579,654.8341627054131,0,0,383.8524390024214,0,269.5612863799487,839.9369396855471,966.480443066781,
    456.4305442695618,0,0,0,0,0,0,538.5968213696259,301.0507196209036,0,516.9672051691632,0,2116.720132085079,0,0,
    364.9417531794148,505.8478642158904,49.71585485491031,0,875.7494136996373,0,422.1093232122134,483.1205250492866,0,0,0,0,0,440.4817582275221,712.7473672074223,
    0,0,1175.731129470483,0,0,6.343491143219572,326.2223003029995,0,0,592.0759328233112,4.14112933427343,0,1693.51092794478,0,0,
    550.7282619328104,0,279.5800833635119,0,0,508.6280344172346,274.568626857958,398.3613296584771,0,362.1602017482217,91.00872626640547,0,26.37600128718807,0,757.3732104433221,
    521.56179378771,144.7735365136191,1100.582851334681,380.5596623712664,0,0,0,0,0,1625.133397696552,66.95680018506519,0,374.0237238963707,1604.92852558617,1140.42884128899,
    0,401.8713775162385,101.5748325918194,4.017421968766477,0,0,0,424.277511022607,0,716.8595298642151,0,0,551.4510758324921,0,0,
    0,4.180334882370253,203.4497540855666,356.7595246265323,0,0,0,888.0596987199661,335.0385551121453,162.173564326569,0,0,87.18345239717291,0,4100.235942741306,
    1240.187112745012,0,270.4710139695795,0,391.3778895655582,0,719.8294609981504,183.310516132936,0,3967.687695871331,527.1016605100402,0,0,173.6080789825696,0,
    9.159861272840839,0,1588.813927571137,0,0,0,0,29.6467403011354,428.1113819448542,197.5793637359868,0,0,790.2220066049952,221.7565058656193,230.6045682097962,
    0,773.1440024699824,0,360.5733264333106,82.81302175650032,98.80607332065935,1859.913651139648,419.3181796201226,6.728421166285229,1797.432442885617,293.0153611569751,474.8601643353819,43.90132772522166,92.73220714306903,0,
    1428.07979808217,0,0,148.4617164181293,0,2983.917996468685,0,19.43350926173421,369.2824085000087,2102.7592437459,328.7896903726656,177.6227646357466,0,370.987751569117,1315.79485664931,
    150.7978242293471,139.1084214258954,185.0844374268527,0,0,370.1045567443864,0,833.1925505915531,88.55448443665929,0,0,0,127.4220265723261,47.80265051669485,0,
    479.4181421371447,0,281.5184696296407,879.2861010997292,0,0,0,195.2892217993547,96.15375078030024,0.3010481265208884,623.7137856424222,165.142213124397,286.3096409000246,909.7818129518311,231.5332093716811,
    0,297.3607462296208,0,88.47278797249362,0,0,460.2791169005494,879.5490008866994,14.90839753769235,141.1267970088761,296.570212477952,708.4443451875763,0,831.5799268538891,96.04516140957479,
    0,0,4.258635817862444,146.515484339549,2565.55405869069,286.6068810596273,0,8.163935601678475,39.81787631149739,2868.215053748437,0,0,0,171.0954318003153,2092.402857532144,
    0,-0,0,0,141.334394842977,0,45.06926595535835,2095.154294035495,222.2340036996875,1560.838736702703,0,238.2377191706084,13.20119188387856,0,0,
    0,0,235.317391615093,0,5.178426678270466,500.138930589119,0,517.1913647430897,247.79402800643,39.92963017325066,2750.35767078554,348.638499302222,0,0,685.3914602409797,
    0,0,496.0004593324446,0,384.9497740474512,437.7309570035949,1427.276572693093,0,0,0,2.663931160723406,3.927970646486636,922.9655196936156,14.67668921775386,1207.388409045948,
    0,514.4722816839188,0,0,154.0526064457605,0,0,243.7167138147433,0,217.8156277060863,0,0,516.2024607828251,0,0,
    7.941835901968208,1479.713680520468,0.633886221903501,131.7256554944386,0,114.5944880682421,0,11.77563732584612,0,794.2081433816187,0,220.9724101158841,0,71.19343974655834,0,
    455.1492571583086,561.5170470550451,0,744.9924720896342,743.3904418346261,0,486.9617830275944,0,0,109.7523918781267,235.1892106461534,804.2824357595298,252.4028226664462,399.8114115274718,76.4691331729224,
    424.7854978720217,0,0,338.4072973177444,219.5942620891569,0,19.81260165722914,425.6729116572579,0,0,0,801.795517618379,0.3020928180076377,67.45582890770814,0,
    972.5534636888488,0,645.1232784538919,254.9005357426872,1.715890854318166,703.3216310043882,2662.62467501002,2991.629304318846,169.5282351320011,0,0,0,127.7404252647524,0,357.6696019834284,
    0,0,0,261.8606883723068,0,104.6852726001587,346.5259307758783,831.8191324207736,0,0,0,530.9128098985983,279.5103593545676,658.2720782024774,561.1058489172541,
    91.95109350313182,0,0,313.1628957753986,0,583.396725683145
)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(47, $i + 1).Value = $rowValues[$i]
}
